$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (D) and Volume(1h) (E) updates ---
$ws.Range("D2").Value = "65.683.38"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").Value = "2.670.71"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'602.40"
$ws.Range("E5").Value = "  -1.38%  "
$ws.Range("D6").Value = "'157.26"
$ws.Range("E6").Value = "  -2.21%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.619"
$ws.Range("E8").Value = "  +4.51%  "
$ws.Range("E9").Value = "  +4.40%  "
$ws.Range("D10").Value = "'0.402"
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("D11").Value = "'5.85"
$ws.Range("E11").Value = "  -2.55%  "
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "'29.51"
$ws.Range("E13").Value = "  -2.48%  "
$ws.Range("E14").Value = "  -4.41%  "
$ws.Range("D15").Value = "3.151.29"
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("D16").Value = "65.527.95"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").Value = "2.634.82"
$ws.Range("E17").Value = "  -1.70%  "
$ws.Range("D18").Value = "'12.93"
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("E19").Value = "  -2.05%  "
$ws.Range("D20").Value = "'7.70"
$ws.Range("E20").Value = "  +2.02%  "
$ws.Range("D21").Value = "'352.18"
$ws.Range("E21").Value = "  -3.07%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "'69.84"
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("E24").Value = "  +3.34%  "
$ws.Range("D25").Value = "'9.77"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("E26").Value = "  -1.57%  "
$ws.Range("D27").Value = "'0.167"
$ws.Range("E27").Value = "  -4.17%  "
$ws.Range("E28").Value = "  -5.64%  "
$ws.Range("E29").Value = "  -1.21%  "
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").Value = "'2.15"
$ws.Range("E31").Value = "  -2.80%  "
$ws.Range("D32").Value = "'531.98"
$ws.Range("E32").Value = "  -1.87%  "
$ws.Range("E33").Value = "  -2.34%  "
$ws.Range("D34").Value = "'6.54"
$ws.Range("E34").Value = "  -1.36%  "
$ws.Range("D35").Value = "'5.50"
$ws.Range("E35").Value = "  +1.36%  "
$ws.Range("E36").Value = "  -2.37%  "
$ws.Range("E37").Value = "  -1.70%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").Value = "'159.04"
$ws.Range("E39").Value = "  -2.59%  "
$ws.Range("E40").Value = "  -3.58%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").Value = "'42.70"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "'165.22"
$ws.Range("E43").Value = "  -3.53%  "
$ws.Range("E44").Value = "  -2.86%  "
$ws.Range("E45").Value = "  -1.38%  "
$ws.Range("D46").Value = "'0.0610"
$ws.Range("E46").Value = "  -1.81%  "
$ws.Range("D47").Value = "'23.13"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("E50").Value = "  +1.72%  "
$ws.Range("D51").Value = "'20.21"
$ws.Range("E51").Value = "  -0.99%  "

# --- Rows 48/49: Mantle and VeChain swapped position, with updated data ---
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.645"
$ws.Range("E48").Value = "  -2.96%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0260"
$ws.Range("E49").Value = "  -3.06%  "
